# Refresh the cryptos price/volume table cells to match the latest scrape.
# (GitHub Actions "Updated cryptos list" run.)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "26.129.39"
$ws.Cells.Item(2, 5).Value = "  -1.46%  "
$ws.Cells.Item(3, 4).Value = "1.657.21"
$ws.Cells.Item(3, 5).Value = "  -1.19%  "
$ws.Cells.Item(4, 5).Value = "  +0.24%  "
$ws.Cells.Item(5, 4).Value = "'216.42"
$ws.Cells.Item(5, 4).Style = "Normal"
$ws.Cells.Item(5, 5).Value = "  -1.53%  "
$ws.Cells.Item(6, 4).Value = "'0.5152"
$ws.Cells.Item(6, 4).Style = "Normal"
$ws.Cells.Item(6, 5).Value = "  -3.16%  "
$ws.Cells.Item(7, 5).Value = "  +0.23%  "
$ws.Cells.Item(8, 4).Value = "'0.2633"
$ws.Cells.Item(8, 4).Style = "Normal"
$ws.Cells.Item(8, 5).Value = "  -2.23%  "
$ws.Cells.Item(9, 4).Value = "'0.06263"
$ws.Cells.Item(9, 4).Style = "Normal"
$ws.Cells.Item(9, 5).Value = "  -2.21%  "
$ws.Cells.Item(10, 4).Value = "'20.71"
$ws.Cells.Item(10, 4).Style = "Normal"
$ws.Cells.Item(10, 5).Value = "  -4.95%  "
$ws.Cells.Item(11, 4).Value = "'0.07721"
$ws.Cells.Item(11, 4).Style = "Normal"
$ws.Cells.Item(11, 5).Value = "  -1.00%  "
$ws.Cells.Item(12, 2).Value = "WrappedEther"
$ws.Cells.Item(12, 3).Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Cells.Item(12, 4).Value = "1.656.69"
$ws.Cells.Item(12, 5).Value = "  -1.81%  "
$ws.Cells.Item(13, 2).Value = "Polkadot"
$ws.Cells.Item(13, 3).Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Cells.Item(13, 4).Value = "'4.428"
$ws.Cells.Item(13, 4).Style = "Normal"
$ws.Cells.Item(13, 5).Value = "  -1.74%  "
$ws.Cells.Item(14, 4).Value = "1.883.98"
$ws.Cells.Item(14, 5).Value = "  -1.23%  "
$ws.Cells.Item(15, 4).Value = "'0.5408"
$ws.Cells.Item(15, 4).Style = "Normal"
$ws.Cells.Item(15, 5).Value = "  -3.19%  "
$ws.Cells.Item(16, 4).Value = "0.0₅8091"
$ws.Cells.Item(16, 5).Value = "  -2.96%  "
$ws.Cells.Item(17, 4).Value = "'64.73"
$ws.Cells.Item(17, 4).Style = "Normal"
$ws.Cells.Item(17, 5).Value = "  -1.56%  "
$ws.Cells.Item(18, 4).Value = "26.159.51"
$ws.Cells.Item(18, 5).Value = "  -1.44%  "
$ws.Cells.Item(20, 4).Value = "'4.601"
$ws.Cells.Item(20, 4).Style = "Normal"
$ws.Cells.Item(20, 5).Value = "  -3.91%  "
$ws.Cells.Item(21, 4).Value = "'191.31"
$ws.Cells.Item(21, 4).Style = "Normal"
$ws.Cells.Item(21, 5).Value = "  -0.74%  "
$ws.Cells.Item(22, 5).Value = "  -2.57%  "
$ws.Cells.Item(23, 4).Value = "'5.997"
$ws.Cells.Item(23, 4).Style = "Normal"
$ws.Cells.Item(23, 5).Value = "  -5.18%  "
$ws.Cells.Item(24, 5).Value = "  +0.27%  "
$ws.Cells.Item(25, 4).Value = "'139.63"
$ws.Cells.Item(25, 4).Style = "Normal"
$ws.Cells.Item(25, 5).Value = "  +0.39%  "
$ws.Cells.Item(26, 4).Value = "'0.1220"
$ws.Cells.Item(26, 4).Style = "Normal"
$ws.Cells.Item(26, 5).Value = "  -4.32%  "
$ws.Cells.Item(27, 4).Value = "'7.191"
$ws.Cells.Item(27, 4).Style = "Normal"
$ws.Cells.Item(27, 5).Value = "  -3.16%  "
$ws.Cells.Item(28, 4).Value = "'16.06"
$ws.Cells.Item(28, 4).Style = "Normal"
$ws.Cells.Item(28, 5).Value = "  -1.41%  "
$ws.Cells.Item(29, 4).Value = "'1.426"
$ws.Cells.Item(29, 4).Style = "Normal"
$ws.Cells.Item(29, 5).Value = "  -0.46%  "
$ws.Cells.Item(30, 4).Value = "'0.05957"
$ws.Cells.Item(30, 4).Style = "Normal"
$ws.Cells.Item(30, 5).Value = "  -5.78%  "
$ws.Cells.Item(31, 4).Value = "'1.269"
$ws.Cells.Item(31, 4).Style = "Normal"
$ws.Cells.Item(32, 4).Value = "'3.553"
$ws.Cells.Item(32, 4).Style = "Normal"
$ws.Cells.Item(32, 5).Value = "  -1.45%  "
$ws.Cells.Item(33, 4).Value = "'3.245"
$ws.Cells.Item(33, 4).Style = "Normal"
$ws.Cells.Item(33, 5).Value = "  -5.71%  "
$ws.Cells.Item(34, 4).Value = "'1.599"
$ws.Cells.Item(34, 4).Style = "Normal"
$ws.Cells.Item(34, 5).Value = "  -5.47%  "
$ws.Cells.Item(35, 4).Value = "'0.9616"
$ws.Cells.Item(35, 4).Style = "Normal"
$ws.Cells.Item(35, 5).Value = "  -5.22%  "
$ws.Cells.Item(36, 4).Value = "'2.426"
$ws.Cells.Item(36, 4).Style = "Normal"
$ws.Cells.Item(36, 5).Value = "  +0.07%  "
$ws.Cells.Item(37, 5).Value = "  -0.53%  "
$ws.Cells.Item(38, 4).Value = "'0.5667"
$ws.Cells.Item(38, 4).Style = "Normal"
$ws.Cells.Item(38, 5).Value = "  -7.81%  "
$ws.Cells.Item(39, 4).Value = "'0.01587"
$ws.Cells.Item(39, 4).Style = "Normal"
$ws.Cells.Item(39, 5).Value = "  -2.82%  "
$ws.Cells.Item(40, 4).Value = "'5.943"
$ws.Cells.Item(40, 4).Style = "Normal"
$ws.Cells.Item(40, 5).Value = "  -2.55%  "
$ws.Cells.Item(41, 4).Value = "'0.8538"
$ws.Cells.Item(41, 4).Style = "Normal"
$ws.Cells.Item(41, 5).Value = "  -0.95%  "
$ws.Cells.Item(43, 4).Value = "1.005.81"
$ws.Cells.Item(43, 5).Value = "  -8.14%  "
$ws.Cells.Item(44, 4).Value = "'100.57"
$ws.Cells.Item(44, 4).Style = "Normal"
$ws.Cells.Item(44, 5).Value = "  -0.06%  "
$ws.Cells.Item(45, 4).Value = "1.798.65"
$ws.Cells.Item(45, 5).Value = "  -1.30%  "
$ws.Cells.Item(46, 4).Value = "'56.66"
$ws.Cells.Item(46, 4).Style = "Normal"
$ws.Cells.Item(46, 5).Value = "  -3.30%  "
$ws.Cells.Item(47, 2).Value = "Frax"
$ws.Cells.Item(47, 3).Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Cells.Item(47, 4).Value = "'1.004"
$ws.Cells.Item(47, 4).Style = "Normal"
$ws.Cells.Item(47, 5).Value = "  +0.52%  "
$ws.Cells.Item(48, 2).Value = "BabyDogeCoin"
$ws.Cells.Item(48, 3).Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Cells.Item(48, 4).Value = "0.0₈108"
$ws.Cells.Item(48, 5).Value = "  -4.22%  "
$ws.Cells.Item(49, 4).Value = "'8.001"
$ws.Cells.Item(49, 4).Style = "Normal"
$ws.Cells.Item(49, 5).Value = "  -2.33%  "
$ws.Cells.Item(50, 5).Value = "  -0.53%  "
$ws.Cells.Item(51, 2).Value = "Mantle"
$ws.Cells.Item(51, 3).Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Cells.Item(51, 4).Value = "'0.4181"
$ws.Cells.Item(51, 4).Style = "Normal"
$ws.Cells.Item(51, 5).Value = "  -1.25%  "
